# "Revert to 2.1.1 files" — fix the broken EPA "mortality risk valuation"
# citation on the About sheet: the old epa.gov short-link (which 404s) is
# replaced by the long-form yosemite.epa.gov URL, and the now-pointless
# live hyperlink object on that cell is removed (the cell keeps its
# Hyperlink-style formatting, it just no longer navigates anywhere).
#
# Also turn off iterative calculation (the workbook no longer needs it),
# matching the source file's calcPr losing its iterate/iterateDelta flags.
$excel.Iteration = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$cell = $ws.Range("B6")

# Remove the live hyperlink attached to B6 (it pointed at the old broken
# epa.gov URL). This also drops the now-unused external relationship.
if ($cell.Hyperlinks.Count -gt 0) {
    $cell.Hyperlinks.Delete()
}

# Swap the cell text for the replacement citation URL. The cell's existing
# "Hyperlink" style (s="3" — blue/underlined font) is left untouched.
$cell.Value = "http://yosemite.epa.gov/EE%5Cepa%5Ceed.nsf/webpages/MortalityRiskValuation.html#whatvalue"
